# Updated cryptos list on Fri Sep  1 23:00:43 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.013.72'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').Value = '1.639.18'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.01'
$ws.Range('E5').Value = '  -1.51%  '
$ws.Range('E6').Value = '  -2.38%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.06439'
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2571'
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.43'
$ws.Range('E10').Value = '  -2.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07723'
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').Value = '1.651.59'
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.238'
$ws.Range('E13').Value = '  -1.37%  '
$ws.Range('D14').Value = '1.865.70'
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5440'
$ws.Range('E15').Value = '  -1.71%  '
$ws.Range('D16').Value = '0.0₅7908'
$ws.Range('E16').Value = '  -1.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.64'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').Value = '26.011.53'
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '203.48'
$ws.Range('E20').Value = '  -4.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.287'
$ws.Range('E21').Value = '  -2.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.972'
$ws.Range('E22').Value = '  -1.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.943'
$ws.Range('E23').Value = '  +0.56%  '
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.927'
$ws.Range('E25').Value = '  +9.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.47'
$ws.Range('E26').Value = '  -1.37%  '
$ws.Range('E27').Value = '  -1.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.72'
$ws.Range('E28').Value = '  -0.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.720'
$ws.Range('E29').Value = '  -3.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05049'
$ws.Range('E30').Value = '  -4.42%  '
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.249'
$ws.Range('E32').Value = '  -3.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.192'
$ws.Range('E33').Value = '  -1.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.540'
$ws.Range('E34').Value = '  -2.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.341'
$ws.Range('E35').Value = '  -0.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.640'
$ws.Range('E36').Value = '  -4.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.8884'
$ws.Range('E37').Value = '  -4.07%  '
$ws.Range('E38').Value = '  -1.70%  '
$ws.Range('D39').Value = '1.144.46'
$ws.Range('E39').Value = '  -1.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01571'
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.563'
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.655'
$ws.Range('E43').Value = '  -0.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8078'
$ws.Range('E44').Value = '  -3.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '99.89'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').Value = '1.777.75'
$ws.Range('D47').Value = '0.0₈112'
$ws.Range('E47').Value = '  +3.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4529'
$ws.Range('E48').Value = '  +0.50%  '
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '54.84'
$ws.Range('E50').Value = '  -2.16%  '
$ws.Range('E51').Value = '  -0.98%  '
